$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.6826468759216482
$ws.Cells.Item(2, 3).Value = 0.1574469820557738
$ws.Cells.Item(2, 4).Value = 0.2803466517295021
$ws.Cells.Item(2, 6).Value = 1.19783137497874
$ws.Cells.Item(2, 7).Value = 0.002426350159526877
$ws.Cells.Item(2, 9).Value = 0.5013199218718007
$ws.Cells.Item(2, 10).Value = 0.269323055743115
$ws.Cells.Item(2, 13).Value = 0.3900585773189533
$ws.Cells.Item(2, 14).Value = 1.161262285041154
$ws.Cells.Item(2, 15).Value = 2.562930756647688

$ws.Cells.Item(3, 2).Value = 0.607081469216638
$ws.Cells.Item(3, 3).Value = 0.1377339995761986
$ws.Cells.Item(3, 4).Value = 0.277881639103569
$ws.Cells.Item(3, 6).Value = 1.195546572664348
$ws.Cells.Item(3, 7).Value = 0.002429068525775315
$ws.Cells.Item(3, 9).Value = 0.5061891078789955
$ws.Cells.Item(3, 10).Value = 0.2694069276888129
$ws.Cells.Item(3, 13).Value = 0.3663935105153371
$ws.Cells.Item(3, 14).Value = 1.170863466632966
$ws.Cells.Item(3, 15).Value = 2.562067369321369

$ws.Cells.Item(4, 2).Value = 0.5606740724697943
$ws.Cells.Item(4, 3).Value = 0.1256009657145682
$ws.Cells.Item(4, 4).Value = 0.2764740242629671
$ws.Cells.Item(4, 6).Value = 1.194867935849445
$ws.Cells.Item(4, 7).Value = 0.002430827864883849
$ws.Cells.Item(4, 9).Value = 0.5095022789070782
$ws.Cells.Item(4, 10).Value = 0.2696008630980842
$ws.Cells.Item(4, 13).Value = 0.3519866908456919
$ws.Cells.Item(4, 14).Value = 1.177213086768525
$ws.Cells.Item(4, 15).Value = 2.563227364524266

$ws.Cells.Item(5, 2).Value = 0.5417613287737026
$ws.Cells.Item(5, 3).Value = 0.1206495623312662
$ws.Cells.Item(5, 4).Value = 0.2759271251445625
$ws.Cells.Item(5, 6).Value = 1.194773461091209
$ws.Cells.Item(5, 7).Value = 0.002431567572129331
$ws.Cells.Item(5, 9).Value = 0.5109336715441444
$ws.Cells.Item(5, 10).Value = 0.2697157040752955
$ws.Cells.Item(5, 13).Value = 0.346147209424096
$ws.Cells.Item(5, 14).Value = 1.179914970788651
$ws.Cells.Item(5, 15).Value = 2.564124719432044

$ws.Cells.Item(6, 2).Value = 0.5386208341331269
$ws.Cells.Item(6, 3).Value = 0.1198269643539902
$ws.Cells.Item(6, 4).Value = 0.2758379289406605
$ws.Cells.Item(6, 6).Value = 1.194768769351079
$ws.Cells.Item(6, 7).Value = 0.002431691776672407
$ws.Cells.Item(6, 9).Value = 0.5111762577522363
$ws.Cells.Item(6, 10).Value = 0.2697369360311512
$ws.Cells.Item(6, 13).Value = 0.3451794750030572
$ws.Cells.Item(6, 14).Value = 1.180370526366339
$ws.Cells.Item(6, 15).Value = 2.564299362081982

$ws.Cells.Item(7, 2).Value = 0.5604190125531261
$ws.Cells.Item(7, 3).Value = 0.1255342177531986
$ws.Cells.Item(7, 4).Value = 0.2764665403105937
$ws.Cells.Item(7, 6).Value = 1.194865924569896
$ws.Cells.Item(7, 7).Value = 0.002430837748623356
$ws.Cells.Item(7, 9).Value = 0.5095212543178427
$ws.Cells.Item(7, 10).Value = 0.2696022668954399
$ws.Cells.Item(7, 13).Value = 0.3519078099260255
$ws.Cells.Item(7, 14).Value = 1.177249062157941
$ws.Cells.Item(7, 15).Value = 2.563237747649168

$ws.Cells.Item(8, 2).Value = 0.6565946478531544
$ws.Cells.Item(8, 3).Value = 0.1506561434115383
$ws.Cells.Item(8, 4).Value = 0.2794747814293146
$ws.Cells.Item(8, 6).Value = 1.196893229894542
$ws.Cells.Item(8, 7).Value = 0.002427268761394874
$ws.Cells.Item(8, 9).Value = 0.5029316260248393
$ws.Cells.Item(8, 10).Value = 0.2693224072699323
$ws.Cells.Item(8, 13).Value = 0.3818733920224346
$ws.Cells.Item(8, 14).Value = 1.16447850069931
$ws.Cells.Item(8, 15).Value = 2.562282059844733

$ws.Cells.Item(9, 2).Value = 0.8450748816859459
$ws.Cells.Item(9, 3).Value = 0.1996801999998752
$ws.Cells.Item(9, 4).Value = 0.2862112197556428
$ws.Cells.Item(9, 6).Value = 1.206619180720949
$ws.Cells.Item(9, 7).Value = 0.002420982997555746
$ws.Cells.Item(9, 9).Value = 0.4925802005728066
$ws.Cells.Item(9, 10).Value = 0.2699044119608374
$ws.Cells.Item(9, 13).Value = 0.4416055219862045
$ws.Cells.Item(9, 14).Value = 1.143037736640437
$ws.Cells.Item(9, 15).Value = 2.573837890717698

$ws.Cells.Item(10, 2).Value = 0.9834365435665973
$ws.Cells.Item(10, 3).Value = 0.2355440368969539
$ws.Cells.Item(10, 4).Value = 0.2916675776696991
$ws.Cells.Item(10, 6).Value = 1.21727875066145
$ws.Cells.Item(10, 7).Value = 0.002416795183429484
$ws.Cells.Item(10, 9).Value = 0.4865482527883032
$ws.Cells.Item(10, 10).Value = 0.2710227460409342
$ws.Cells.Item(10, 13).Value = 0.4860715077379822
$ws.Cells.Item(10, 14).Value = 1.129476545648821
$ws.Cells.Item(10, 15).Value = 2.590548055631785

$ws.Cells.Item(11, 2).Value = 1.046347925895873
$ws.Cells.Item(11, 3).Value = 0.2518245026421653
$ws.Cells.Item(11, 4).Value = 0.2942592726429467
$ws.Cells.Item(11, 6).Value = 1.222893054531411
$ws.Cells.Item(11, 7).Value = 0.002414982560581766
$ws.Cells.Item(11, 9).Value = 0.4841470888273562
$ws.Cells.Item(11, 10).Value = 0.2716818199838826
$ws.Cells.Item(11, 13).Value = 0.506424437423469
$ws.Cells.Item(11, 14).Value = 1.123782138714688
$ws.Cells.Item(11, 15).Value = 2.599942288388036

$ws.Cells.Item(12, 2).Value = 1.070165531535963
$ws.Cells.Item(12, 3).Value = 0.2579843756084585
$ws.Cells.Item(12, 4).Value = 0.2952563646745006
$ws.Cells.Item(12, 6).Value = 1.225129185639148
$ws.Cells.Item(12, 7).Value = 0.002414309389642932
$ws.Cells.Item(12, 9).Value = 0.4832872316747938
$ws.Cells.Item(12, 10).Value = 0.2719530304988353
$ws.Cells.Item(12, 13).Value = 0.51414928885373
$ws.Cells.Item(12, 14).Value = 1.12169401125692
$ws.Cells.Item(12, 15).Value = 2.603757920455337

$ws.Cells.Item(13, 2).Value = 1.065036251608433
$ws.Cells.Item(13, 3).Value = 0.2566579717287993
$ws.Cells.Item(13, 4).Value = 0.2950409273026224
$ws.Cells.Item(13, 6).Value = 1.22464269609462
$ws.Cells.Item(13, 7).Value = 0.002414453781633631
$ws.Cells.Item(13, 9).Value = 0.4834702178387609
$ws.Cells.Item(13, 10).Value = 0.2718936580650393
$ws.Cells.Item(13, 13).Value = 0.5124848250456324
$ws.Cells.Item(13, 14).Value = 1.122140693201459
$ws.Cells.Item(13, 15).Value = 2.602924663337859

$ws.Cells.Item(14, 2).Value = 1.048307533424349
$ws.Cells.Item(14, 3).Value = 0.2523313845905193
$ws.Cells.Item(14, 4).Value = 0.2943409903589469
$ws.Cells.Item(14, 6).Value = 1.223074815040434
$ws.Cells.Item(14, 7).Value = 0.002414926913636796
$ws.Cells.Item(14, 9).Value = 0.484075356909468
$ws.Cells.Item(14, 10).Value = 0.2717036989818453
$ws.Cells.Item(14, 13).Value = 0.5070596141261774
$ws.Cells.Item(14, 14).Value = 1.123608980137057
$ws.Cells.Item(14, 15).Value = 2.600251024634701

$ws.Cells.Item(15, 2).Value = 1.038059948104717
$ws.Cells.Item(15, 3).Value = 0.2496805415970016
$ws.Cells.Item(15, 4).Value = 0.2939142974569506
$ws.Cells.Item(15, 6).Value = 1.222128785334561
$ws.Cells.Item(15, 7).Value = 0.002415218441586859
$ws.Cells.Item(15, 9).Value = 0.484452460384059
$ws.Cells.Item(15, 10).Value = 0.2715901613284615
$ws.Cells.Item(15, 13).Value = 0.5037388023370113
$ws.Cells.Item(15, 14).Value = 1.124517231751277
$ws.Cells.Item(15, 15).Value = 2.598646987600887

$ws.Cells.Item(16, 2).Value = 0.979324425473294
$ws.Cells.Item(16, 3).Value = 0.2344793592552037
$ws.Cells.Item(16, 4).Value = 0.2915004020737797
$ws.Cells.Item(16, 6).Value = 1.216927247536262
$ws.Cells.Item(16, 7).Value = 0.002416915497442726
$ws.Cells.Item(16, 9).Value = 0.486712082538677
$ws.Cells.Item(16, 10).Value = 0.2709827004849217
$ws.Cells.Item(16, 13).Value = 0.4847438835108733
$ws.Cells.Item(16, 14).Value = 1.129858236742237
$ws.Cells.Item(16, 15).Value = 2.589970232771918

$ws.Cells.Item(17, 2).Value = 0.9432835021231654
$ws.Cells.Item(17, 3).Value = 0.2251449799428826
$ws.Cells.Item(17, 4).Value = 0.2900475601790333
$ws.Cells.Item(17, 6).Value = 1.213932318545105
$ws.Cells.Item(17, 7).Value = 0.002417980216562781
$ws.Cells.Item(17, 9).Value = 0.4881861702620327
$ws.Cells.Item(17, 10).Value = 0.2706485588601026
$ws.Cells.Item(17, 13).Value = 0.4731229122491598
$ws.Cells.Item(17, 14).Value = 1.133256314076505
$ws.Cells.Item(17, 15).Value = 2.58510679960645

$ws.Cells.Item(18, 2).Value = 0.9225509641654526
$ws.Cells.Item(18, 3).Value = 0.2197728915020321
$ws.Cells.Item(18, 4).Value = 0.2892222421702968
$ws.Cells.Item(18, 6).Value = 1.212281735266757
$ws.Cells.Item(18, 7).Value = 0.002418601319137951
$ws.Cells.Item(18, 9).Value = 0.4890662854172625
$ws.Cells.Item(18, 10).Value = 0.2704705190732
$ws.Cells.Item(18, 13).Value = 0.4664506359437866
$ws.Cells.Item(18, 14).Value = 1.13525547597164
$ws.Cells.Item(18, 15).Value = 2.582478191459387

$ws.Cells.Item(19, 2).Value = 0.91553084852967
$ws.Cells.Item(19, 3).Value = 0.2179534532934042
$ws.Cells.Item(19, 4).Value = 0.2889445783592635
$ws.Cells.Item(19, 6).Value = 1.211735243316468
$ws.Cells.Item(19, 7).Value = 0.002418813110608797
$ws.Cells.Item(19, 9).Value = 0.4893698142026537
$ws.Cells.Item(19, 10).Value = 0.2704126676155667
$ws.Cells.Item(19, 13).Value = 0.4641935551374914
$ws.Cells.Item(19, 14).Value = 1.135940032347008
$ws.Cells.Item(19, 15).Value = 2.581617152191001

$ws.Cells.Item(20, 2).Value = 0.9471204138403095
$ws.Cells.Item(20, 3).Value = 0.2261389737511195
$ws.Cells.Item(20, 4).Value = 0.2902011503810229
$ws.Cells.Item(20, 6).Value = 1.214243679400496
$ws.Cells.Item(20, 7).Value = 0.002417865975126692
$ws.Cells.Item(20, 9).Value = 0.4880259115818966
$ws.Cells.Item(20, 10).Value = 0.270682664262992
$ws.Cells.Item(20, 13).Value = 0.4743587658974917
$ws.Cells.Item(20, 14).Value = 1.132889959219135
$ws.Cells.Item(20, 15).Value = 2.585607056229861

$ws.Cells.Item(21, 2).Value = 1.053221324577692
$ws.Cells.Item(21, 3).Value = 0.2536023510258474
$ws.Cells.Item(21, 4).Value = 0.2945461540937089
$ws.Cells.Item(21, 6).Value = 1.223532350563104
$ws.Cells.Item(21, 7).Value = 0.002414787584530323
$ws.Cells.Item(21, 9).Value = 0.483896270946147
$ws.Cells.Item(21, 10).Value = 0.2717589073214555
$ws.Cells.Item(21, 13).Value = 0.5086526542005814
$ws.Cells.Item(21, 14).Value = 1.123175857417664
$ws.Cells.Item(21, 15).Value = 2.601029324974292

$ws.Cells.Item(22, 2).Value = 1.122531361066422
$ws.Cells.Item(22, 3).Value = 0.2715209013419155
$ws.Cells.Item(22, 4).Value = 0.2974771834054764
$ws.Cells.Item(22, 6).Value = 1.230244911299636
$ws.Cells.Item(22, 7).Value = 0.002412852765370173
$ws.Cells.Item(22, 9).Value = 0.4814853746232544
$ws.Cells.Item(22, 10).Value = 0.2725883980773602
$ws.Cells.Item(22, 13).Value = 0.5311683398481506
$ws.Cells.Item(22, 14).Value = 1.117224761537152
$ws.Cells.Item(22, 15).Value = 2.612614098175328

$ws.Cells.Item(23, 2).Value = 1.085542726275605
$ws.Cells.Item(23, 3).Value = 0.261960300980121
$ws.Cells.Item(23, 4).Value = 0.2959045095777952
$ws.Cells.Item(23, 6).Value = 1.226603533413865
$ws.Cells.Item(23, 7).Value = 0.002413878381698242
$ws.Cells.Item(23, 9).Value = 0.4827457188750159
$ws.Cells.Item(23, 10).Value = 0.2721341397666563
$ws.Cells.Item(23, 13).Value = 0.5191420196333496
$ws.Cells.Item(23, 14).Value = 1.120364599855996
$ws.Cells.Item(23, 15).Value = 2.606293190787767

$ws.Cells.Item(24, 2).Value = 0.9453857837448822
$ws.Cells.Item(24, 3).Value = 0.2256896066666343
$ws.Cells.Item(24, 4).Value = 0.2901316812827304
$ws.Cells.Item(24, 6).Value = 1.214102691250545
$ws.Cells.Item(24, 7).Value = 0.0024179175958605
$ws.Cells.Item(24, 9).Value = 0.4880982628596477
$ws.Cells.Item(24, 10).Value = 0.2706672014080951
$ws.Cells.Item(24, 13).Value = 0.4738000091041528
$ws.Cells.Item(24, 14).Value = 1.133055446322267
$ws.Cells.Item(24, 15).Value = 2.585380368660992

$ws.Cells.Item(25, 2).Value = 0.7941031378873618
$ws.Cells.Item(25, 3).Value = 0.1864443967094473
$ws.Cells.Item(25, 4).Value = 0.2842995426528319
$ws.Cells.Item(25, 6).Value = 1.203371571807992
$ws.Cells.Item(25, 7).Value = 0.002422607580504005
$ws.Cells.Item(25, 9).Value = 0.4951046909548396
$ws.Cells.Item(25, 10).Value = 0.2696257626814997
$ws.Cells.Item(25, 13).Value = 0.4253436571821112
$ws.Cells.Item(25, 14).Value = 1.148452886964876
$ws.Cells.Item(25, 15).Value = 2.569270565641773
